# "Update About Me slide." -- edits to slide 2 (the "About Me" slide):
#  - reflow/resize the bullet text box and split "15+ Years..." into
#    "20+ " + "Years working with SQL Server"
#  - center the title and give it an explicit position/size
#  - reposition the headshot, icons, and caption text boxes
#  - widen/move the footer textbox and change the URL it shows

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Shape 1: "Content Placeholder 2" (bullet list starting "15+ Years ...") ---
$sh1 = $s.Shapes.Item(1)
$sh1.Left = 580.1143798828125
$sh1.Top = 120.32921600341797
$sh1.Width = 301.7142639160156
$sh1.Height = 288.64630126953125

# Split the leading "15+ " into its own run reading "20+ " and keep the
# remainder of the sentence in the original run.
$tr1 = $sh1.TextFrame.TextRange
$lead = $tr1.Characters(1, 4)
$lead.Text = "20+ "

# --- Shape 2: "Title 1" (Brian Hansen) ---
$sh2 = $s.Shapes.Item(2)
$sh2.Left = 152.54244995117188
$sh2.Top = 150.23118591308594
$sh2.Width = 294.5432434082031
$sh2.Height = 56.69291687011719
$sh2.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- Shape 3: "Content Placeholder 5" (headshot picture) ---
$sh3 = $s.Shapes.Item(3)
$sh3.Left = 336.62152099609375
$sh3.Top = 254.57615661621094

# --- Shape 4: "Picture 6" (small icon) ---
$sh4 = $s.Shapes.Item(4)
$sh4.Left = 45.24110412597656
$sh4.Top = 301.2997741699219

# --- Shape 5: "Picture 7" (small icon) ---
$sh5 = $s.Shapes.Item(5)
$sh5.Left = 45.24110412597656
$sh5.Top = 255.40158081054688

# --- Shape 6: "Picture 12" (logo image) ---
$sh6 = $s.Shapes.Item(6)
$sh6.Left = 214.69134521484375
$sh6.Top = 75.39614868164062

# --- Shape 7: "Content Placeholder 2" (@tf3604.com) ---
$sh7 = $s.Shapes.Item(7)
$sh7.Left = 87.44763946533203
$sh7.Top = 294.3934020996094

# --- Shape 8: "Content Placeholder 2" (brian@tf3604.com) ---
$sh8 = $s.Shapes.Item(8)
$sh8.Left = 87.44772338867188
$sh8.Top = 258.3822326660156

# --- Shape 9: "Content Placeholder 2" (children.org) ---
$sh9 = $s.Shapes.Item(9)
$sh9.Left = 336.62152099609375
$sh9.Top = 331.3711242675781

# --- Shape 10: "Content Placeholder 2" (footer URL) ---
$sh10 = $s.Shapes.Item(10)
$sh10.Left = 45.24110412597656
$sh10.Top = 419.62835693359375
$sh10.Width = 821.5017700195312
$sh10.Height = 33.598899841308594
$sh10.TextFrame.TextRange.Text = "www.tf3604.com/poshadmin"
